# Submission tracker update: add newly submitted runs (rows 14-16) and
# fill in the "By" info that was missing for the existing last row (13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

function Set-TextValue($Cell, $Text) {
    # Plain assignment auto-types "TRUE"/"FALSE"-looking strings into real
    # booleans (like typing them straight into Excel). The source data
    # keeps them as plain text, so route those through a formula + paste
    # values round-trip, which preserves the literal text.
    if ($Text -eq "TRUE" -or $Text -eq "FALSE") {
        $Cell.Formula = '="' + $Text + '"'
        $Cell.Copy()
        $Cell.PasteSpecial(-4163)
    } else {
        $Cell.Value = $Text
    }
}

# --- Grow the table by 3 rows (13 -> 16) -----------------------------------
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# --- Row 13: fill in the previously-empty "Name" cell -----------------------
Set-TextValue $ws.Range("B13") "221123_xgb_reduced_3"

# --- Row 14 -----------------------------------------------------------------
$ws.Range("A14").Value = 44888
$ws.Range("A14").NumberFormat = "m/d/yy"
Set-TextValue $ws.Range("B14") "221123_cat"
Set-TextValue $ws.Range("C14") "test_JM_MS"
Set-TextValue $ws.Range("D14") "TRUE"
Set-TextValue $ws.Range("E14") "Maria"

# --- Row 15 -----------------------------------------------------------------
$ws.Range("A15").Value = 44889
$ws.Range("A15").NumberFormat = "m/d/yy"
Set-TextValue $ws.Range("B15") "221123_cat_monkey"
Set-TextValue $ws.Range("C15") "test_bis_JM_MS"
Set-TextValue $ws.Range("D15") "TRUE"
Set-TextValue $ws.Range("E15") "Maria"

# --- Row 16 (no submission, just a placeholder "-") --------------------------
$ws.Range("A16").NumberFormat = "m/d/yy"
Set-TextValue $ws.Range("A16") "-"
Set-TextValue $ws.Range("B16") "221122_ET"
Set-TextValue $ws.Range("D16") "FALSE"

# --- View bookkeeping: zoom level + active selection moved to B15 -----------
$excel.ActiveWindow.Zoom = 177
$ws.Range("B15").Select()

Write-Output "done"
